# Fix formatting issues introduced when scraping data:
#  1) A handful of "Razon social"/"Nombre Fantasia" entries used commas to
#     separate multiple people/abbreviations; these are normalized to periods
#     (and a redundant "S.H." abbreviation is tidied to "SH").
#  2) The "Importe" column was scraped using Spanish/Argentine number
#     formatting (e.g. "27.500,00", thousands separator ".", decimal ",").
#     These text values are normalized to a plain numeric-looking string
#     (e.g. "27500.00"): the "." thousands separators are removed and the
#     "," decimal separator becomes a ".".
#
# Both columns store their data as plain text (shared strings), not real
# numbers, so a plain ".Value = '27500.00'" assignment would make Excel's
# auto-detection silently coerce the text into a real floating point number
# (losing the trailing ".00" and introducing binary rounding noise). To keep
# these as literal text we stage the new text in a scratch cell formatted as
# Text ("@") and PasteSpecial(xlPasteValues) it into the destination - this
# carries over the literal text, not a parsed number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# Scratch cell, well outside the used range, used to force literal-text
# round-trips through copy / paste-special-values.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count
$scratchRowIdx = $lastRow + 10
$scratch = $ws.Cells.Item($scratchRowIdx, 1)
$scratch.NumberFormat = "@"

function Set-LiteralText($cell, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

# --- 1) Specific name/text corrections (columns E and F) ---
$nameFixes = @{
    "RICCOTTI, MARIANA EDITH" = "RICCOTTI. MARIANA EDITH";
    "ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN";
    "SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H." = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH";
    "ODIARD, OSCAR HERNAN" = "ODIARD. OSCAR HERNAN";
    "MORERA, SALVADOR" = "MORERA. SALVADOR";
}

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in 5, 6) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value()
        if ($val -ne $null -and $nameFixes.ContainsKey($val)) {
            $cell.Value = $nameFixes[$val]
        }
    }
}

# --- 2) Reformat "Importe" (column H) numbers from "1.234,56" to "1234.56" ---
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $val = $cell.Value()
    if ($val -ne $null -and $val -match "^\d{1,3}(\.\d{3})*,\d{2}$") {
        $fixed = $val -replace "\.", ""
        $fixed = $fixed -replace ",", "."
        Set-LiteralText $cell $fixed
    }
}

# Clean up the scratch row entirely so it leaves no trace in the sheet.
$excel.CutCopyMode = $false
$scratch.ClearContents()
$ws.Rows.Item($scratchRowIdx).Delete()
